$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped from
# 45175 (2023-09-06) to 45177 (2023-09-08) for every data row (rows 2-499).
$oldValue = 45175
$newValue = 45177

for ($r = 2; $r -le 499; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
